$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns (I and J) before the existing "On Hit Behavior" /
# "Reload Speed Mult" columns. This shifts the old I,J -> K,L and the new
# I,J inherit the neighboring (column H) number formatting/style.
$ws.Columns("I:J").Insert()

# Match the new columns' display width to column H's width.
$ws.Columns("I:J").ColumnWidth = 21.7

# Headers for the two new columns
$ws.Range("I1").Value = "Base DMG Pistol Mult"
$ws.Range("J1").Value = "Out of Range Pistol Mult"

# Fill in the Pistol multiplier values for every data row (2-37)
$ws.Range("I2:I37").Value = 0.95
$ws.Range("J2:J37").Value = 0.7

# Re-enter the DMG (Calculated) formula across the full column so Excel
# consolidates it into a single shared formula, matching authoring flow
# (e.g. fill-down) rather than 36 independent formula instances.
$ws.Range("D2:D37").Formula = "=_xlfn.FLOOR.MATH(PRODUCT(E2,F2))"

# Update the active selection, as recorded in the saved view state.
[void]$ws.Range("J38").Select()
